$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: AMP Limited (ASX:AMP) - refreshed financial metrics
# Row 3: Generation Development Group Limited (ASX:GDG) - refreshed financial metrics, historical_growth_revenue_last_5_years (D3) cleared
# Row 4/5: company names swapped, each with refreshed financial metrics

# Clear D3 (historical_growth_revenue_last_5_years no longer available for this company)
$ws.Range("D3").ClearContents()


# Row 2
$ws.Range("D2").Value = 0.044205
$ws.Range("E2").Value = 0.0276
$ws.Range("F2").Value = 0.1
$ws.Range("G2").Value = 0.06604584303370259
$ws.Range("H2").Value = 0.06604584303370259
$ws.Range("I2").Value = 0.0622646521148107
$ws.Range("J2").Value = 0.0501616324017932
$ws.Range("K2").Value = 278.57
$ws.Range("L2").Value = 0.04352588240808737
$ws.Range("M2").Value = 377.6
$ws.Range("N2").Value = 0.04397037589081932
$ws.Range("O2").Value = 1.355494130739132
$ws.Range("P2").Value = 370.63
$ws.Range("Q2").Value = 0.04315874050957194
$ws.Range("R2").Value = 1.330473489607639
$ws.Range("S2").Value = 6.969999999999999
$ws.Range("T2").Value = 0.01845868644067796
$ws.Range("U2").Value = 781.2
$ws.Range("V2").Value = 0.09096837300293446
$ws.Range("W2").Value = 0.1403113016016242
$ws.Range("X2").Value = 0.05520911817519505
$ws.Range("Y2").Value = 0.08510218342642914
$ws.Range("Z2").Value = 5.79605801763423
$ws.Range("AA2").Value = 0.2835153330957075
$ws.Range("AB2").Value = 0.05472240243958535
$ws.Range("AC2").Value = 0.2287929306561221
$ws.Range("AD2").Value = 293.537
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 293.537
$ws.Range("AG2").Value = -487.663
$ws.Range("AH2").Value = 0.03305173650626041
$ws.Range("AI2").Value = 0.1485107539094892
$ws.Range("AJ2").Value = -0.06020577690912905
$ws.Range("AK2").Value = -0.4079711411928185
$ws.Range("AL2").Value = 9.122
$ws.Range("AM2").Value = 9.122
$ws.Range("AN2").Value = 0.6944334043056543
$ws.Range("AO2").Value = 43.68559526419645
$ws.Range("AP2").Value = -1.153685829193281
$ws.Range("AQ2").Value = 43.68559526419645

# Row 3
$ws.Range("F3").Value = 0.139
$ws.Range("G3").Value = 0.5989583333333334
$ws.Range("H3").Value = 0.5989583333333334
$ws.Range("I3").Value = 0.6041666666666666
$ws.Range("J3").Value = 0.6041666666666666
$ws.Range("K3").Value = -1.13
$ws.Range("L3").Value = 0.05885416666666667
$ws.Range("M3").Value = 1.73
$ws.Range("N3").Value = 0.01888646288209607
$ws.Range("O3").Value = -1.530973451327434
$ws.Range("P3").Value = 1.73
$ws.Range("Q3").Value = 0.01888646288209607
$ws.Range("R3").Value = -1.530973451327434
$ws.Range("S3").Value = 0
$ws.Range("U3").Value = 42.9
$ws.Range("V3").Value = 0.4683406113537118
$ws.Range("W3").Value = -0.07793103448275861
$ws.Range("X3").Value = 0.05495059774014913
$ws.Range("Y3").Value = -0.1328816322229077
$ws.Range("Z3").Value = 7.153502235469447
$ws.Range("AA3").Value = 4.321907600596124
$ws.Range("AB3").Value = 0.05481537327732752
$ws.Range("AC3").Value = 4.267092227318797
$ws.Range("AD3").Value = 0.337
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 0.337
$ws.Range("AG3").Value = -42.563
$ws.Range("AH3").Value = 0.003665553585607536
$ws.Range("AI3").Value = 0.02731620329091352
$ws.Range("AJ3").Value = -0.8679772416746537
$ws.Range("AK3").Value = 1.392631613388738
$ws.Range("AL3").Value = 0.07199999999999999
$ws.Range("AM3").Value = 0.07199999999999999
$ws.Range("AN3").Value = -0.02930434782608696
$ws.Range("AO3").Value = -161.1111111111111
$ws.Range("AP3").Value = 3.701130434782608
$ws.Range("AQ3").Value = -161.1111111111111

# Row 4
$ws.Range("B4").Value = 'Medibank Private Limited (ASX:MPL)'
$ws.Range("D4").Value = 0.00291
$ws.Range("E4").Value = 0.02
$ws.Range("F4").Value = 0.0337
$ws.Range("G4").Value = 0.0701833098320728
$ws.Range("H4").Value = 0.0701833098320728
$ws.Range("I4").Value = 0.06691449814126393
$ws.Range("J4").Value = 0.04691334988732758
$ws.Range("K4").Value = 217.5
$ws.Range("L4").Value = 0.04646840148698884
$ws.Range("M4").Value = 299.22
$ws.Range("N4").Value = 0.046787484558973
$ws.Range("O4").Value = 1.375724137931035
$ws.Range("P4").Value = 296.6
$ws.Range("Q4").Value = 0.04637780870326647
$ws.Range("R4").Value = 1.36367816091954
$ws.Range("S4").Value = 2.620000000000005
$ws.Range("T4").Value = 0.008756099191230547
$ws.Range("U4").Value = 601.6
$ws.Range("V4").Value = 0.0940690819820806
$ws.Range("W4").Value = 0.160280029476787
$ws.Range("X4").Value = 0.05520911817519505
$ws.Range("Y4").Value = 0.105070911301592
$ws.Range("Z4").Value = 6.043382827630729
$ws.Range("AA4").Value = 0.2835153330957075
$ws.Range("AB4").Value = 0.05472240243958535
$ws.Range("AC4").Value = 0.2287929306561221
$ws.Range("AD4").Value = 75.40000000000001
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 75.40000000000001
$ws.Range("AG4").Value = -526.2
$ws.Range("AH4").Value = 0.01165252600182361
$ws.Range("AI4").Value = 0.05678993748587784
$ws.Range("AJ4").Value = -0.0896559949566373
$ws.Range("AK4").Value = -0.7246935683790113
$ws.Range("AL4").Value = 2.35
$ws.Range("AM4").Value = 2.35
$ws.Range("AN4").Value = 0.2295281582952816
$ws.Range("AO4").Value = 133.2765957446808
$ws.Range("AP4").Value = -1.601826484018265
$ws.Range("AQ4").Value = 133.2765957446808

# Row 5
$ws.Range("B5").Value = 'nib holdings limited (ASX:NHF)'
$ws.Range("D5").Value = 0.08550000000000001
$ws.Range("E5").Value = 0.0352
$ws.Range("F5").Value = 0.1
$ws.Range("G5").Value = 0.0607925461551734
$ws.Range("H5").Value = 0.0607925461551734
$ws.Range("I5").Value = 0.05573129349514005
$ws.Range("J5").Value = 0.03989052307707378
$ws.Range("K5").Value = 62.2
$ws.Range("L5").Value = 0.03577385402887215
$ws.Range("M5").Value = 76.64999999999999
$ws.Range("N5").Value = 0.03648783738753748
$ws.Range("O5").Value = 1.232315112540193
$ws.Range("P5").Value = 72.3
$ws.Range("Q5").Value = 0.03441709906221736
$ws.Range("R5").Value = 1.162379421221865
$ws.Range("S5").Value = 4.349999999999994
$ws.Range("T5").Value = 0.05675146771037175
$ws.Range("U5").Value = 136.7
$ws.Range("V5").Value = 0.06507354691293378
$ws.Range("W5").Value = 0.1403113016016242
$ws.Range("X5").Value = 0.05813795303571433
$ws.Range("Y5").Value = 0.08217334856590987
$ws.Range("Z5").Value = 5.230746089049338
$ws.Range("AA5").Value = 0.2086571975755361
$ws.Range("AB5").Value = 0.0539390252068687
$ws.Range("AC5").Value = 0.1547181723686674
$ws.Range("AD5").Value = 217.8
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 217.8
$ws.Range("AG5").Value = 81.10000000000002
$ws.Range("AH5").Value = 0.09394004744446841
$ws.Range("AI5").Value = 0.3421838177533386
$ws.Range("AJ5").Value = 0.03717114309285912
$ws.Range("AK5").Value = 0.162264905962385
$ws.Range("AL5").Value = 6.7
$ws.Range("AM5").Value = 6.7
$ws.Range("AN5").Value = 2.060548722800378
$ws.Range("AO5").Value = 14.46268656716418
$ws.Range("AP5").Value = 0.7672658467360456
$ws.Range("AQ5").Value = 14.46268656716418
